$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2:D9").Value = 0

$ws.Range("C4").Value = 0.6518920857689061
$ws.Range("C5").Value = -0.6549810301751726
$ws.Range("C8").Value = 0.6415598863873151
$ws.Range("C9").Value = -0.656754724003772
